{"js": "// Update the title date and the 100 arithmetic-problem cells in the\n// single table. The document body is, in order:\n//   paragraph 0            -> the \"YYYY-MM-DD Weekday\" title\n//   paragraphs 1..100      -> one paragraph per table cell (20 rows x 5 cols),\n//                              each holding a single \"a+b=\" / \"a-b=\" expression.\n// We replace each paragraph's text positionally (not by unique search/replace)\n// because some old expressions (e.g. \"84-12=\") repeat and map to different\n// new values depending on position.\nconst oldValues = [\"2025-11-11 Tuesday\", \"55+16=\", \"39+25=\", \"77-53=\", \"38-11=\", \"31+5=\", \"35-7=\", \"93-83=\", \"16+71=\", \"90-81=\", \"81+0=\", \"89-86=\", \"11+5=\", \"17+77=\", \"43+54=\", \"46+41=\", \"31+46=\", \"0+66=\", \"86-24=\", \"36+4=\", \"32+39=\", \"37+39=\", \"57+13=\", \"17+40=\", \"9+24=\", \"26+25=\", \"45+4=\", \"55-16=\", \"35+44=\", \"97-16=\", \"35-0=\", \"21+69=\", \"9-7=\", \"34-4=\", \"53-39=\", \"45-34=\", \"23+22=\", \"52+25=\", \"97-22=\", \"59+25=\", \"39-26=\", \"82-29=\", \"43-30=\", \"36+50=\", \"85-21=\", \"49+44=\", \"30+55=\", \"75-20=\", \"55-11=\", \"84-12=\", \"15+49=\", \"64-22=\", \"74-10=\", \"2+88=\", \"93-71=\", \"13-6=\", \"34+50=\", \"59+31=\", \"65-41=\", \"19+59=\", \"88-60=\", \"3+34=\", \"0+96=\", \"54+20=\", \"85-28=\", \"20+11=\", \"10+43=\", \"48-6=\", \"7+15=\", \"86-72=\", \"58-12=\", \"54-33=\", \"1+49=\", \"51-4=\", \"25-20=\", \"59-19=\", \"34+46=\", \"6+47=\", \"40+44=\", \"63-25=\", \"88-73=\", \"14-4=\", \"15+75=\", \"44-29=\", \"41-3=\", \"65-25=\", \"95-3=\", \"14+43=\", \"47-21=\", \"2+56=\", \"40+25=\", \"71-62=\", \"90-88=\", \"52-46=\", \"62+27=\", \"76-10=\", \"58-11=\", \"84-12=\", \"10+69=\", \"28-3=\", \"43+36=\"];\nconst newValues = [\"2025-11-12 Wednesday\", \"1+82=\", \"9-4=\", \"20+63=\", \"7+8=\", \"18+33=\", \"90-68=\", \"21+46=\", \"31-2=\", \"74-20=\", \"60-13=\", \"70-22=\", \"71+22=\", \"61-40=\", \"7+44=\", \"87-62=\", \"23+7=\", \"74-7=\", \"2+78=\", \"49+33=\", \"61+35=\", \"84-80=\", \"28+6=\", \"44+26=\", \"45+20=\", \"98-19=\", \"62-60=\", \"13+63=\", \"17+70=\", \"46-28=\", \"55-50=\", \"41-31=\", \"75-60=\", \"82-59=\", \"29+59=\", \"1+84=\", \"38+17=\", \"61+7=\", \"99-12=\", \"39+53=\", \"35-20=\", \"53-6=\", \"41+6=\", \"94-10=\", \"94+3=\", \"20-19=\", \"20+52=\", \"37+31=\", \"78-24=\", \"99-2=\", \"58+0=\", \"43-20=\", \"83-3=\", \"66+8=\", \"72-36=\", \"6+67=\", \"28+25=\", \"58+30=\", \"59-20=\", \"26+15=\", \"5+75=\", \"94-82=\", \"57+42=\", \"41-21=\", \"86-70=\", \"71-63=\", \"69-52=\", \"96-34=\", \"93-30=\", \"17+60=\", \"37-13=\", \"30+52=\", \"51-1=\", \"74+21=\", \"24-2=\", \"51-40=\", \"85-27=\", \"66-20=\", \"23+56=\", \"7+84=\", \"40-0=\", \"95-83=\", \"48-8=\", \"28+59=\", \"38+6=\", \"54-3=\", \"79+1=\", \"79-77=\", \"30+14=\", \"97-65=\", \"66-54=\", \"54-10=\", \"41-7=\", \"85-57=\", \"38-31=\", \"72-69=\", \"59-5=\", \"51+0=\", \"7+47=\", \"68-7=\", \"36-23=\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newValues.length + \" got \" + items.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const current = items[i].text;\n  if (current !== oldValues[i]) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected '\" + oldValues[i] + \"' got '\" + current + \"'\"\n    );\n  }\n  if (current !== newValues[i]) {\n    items[i].insertText(newValues[i], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and the 100 arithmetic-problem cells in the\n# single table (20 rows x 5 columns). Cell values are applied\n# positionally in row-major order (row 1 col 1..5, row 2 col 1..5, ...)\n# because a couple of old expressions (e.g. \"84-12=\") repeat verbatim\n# but map to different new values depending on position, so a global\n# text replace would be ambiguous.\n\n$d = $word.ActiveDocument\n\n$titleOld = \"2025-11-11 Tuesday\"\n$titleNew = \"2025-11-12 Wednesday\"\n\n$oldValues = @(\"55+16=\",\"39+25=\",\"77-53=\",\"38-11=\",\"31+5=\",\"35-7=\",\"93-83=\",\"16+71=\",\"90-81=\",\"81+0=\",\"89-86=\",\"11+5=\",\"17+77=\",\"43+54=\",\"46+41=\",\"31+46=\",\"0+66=\",\"86-24=\",\"36+4=\",\"32+39=\",\"37+39=\",\"57+13=\",\"17+40=\",\"9+24=\",\"26+25=\",\"45+4=\",\"55-16=\",\"35+44=\",\"97-16=\",\"35-0=\",\"21+69=\",\"9-7=\",\"34-4=\",\"53-39=\",\"45-34=\",\"23+22=\",\"52+25=\",\"97-22=\",\"59+25=\",\"39-26=\",\"82-29=\",\"43-30=\",\"36+50=\",\"85-21=\",\"49+44=\",\"30+55=\",\"75-20=\",\"55-11=\",\"84-12=\",\"15+49=\",\"64-22=\",\"74-10=\",\"2+88=\",\"93-71=\",\"13-6=\",\"34+50=\",\"59+31=\",\"65-41=\",\"19+59=\",\"88-60=\",\"3+34=\",\"0+96=\",\"54+20=\",\"85-28=\",\"20+11=\",\"10+43=\",\"48-6=\",\"7+15=\",\"86-72=\",\"58-12=\",\"54-33=\",\"1+49=\",\"51-4=\",\"25-20=\",\"59-19=\",\"34+46=\",\"6+47=\",\"40+44=\",\"63-25=\",\"88-73=\",\"14-4=\",\"15+75=\",\"44-29=\",\"41-3=\",\"65-25=\",\"95-3=\",\"14+43=\",\"47-21=\",\"2+56=\",\"40+25=\",\"71-62=\",\"90-88=\",\"52-46=\",\"62+27=\",\"76-10=\",\"58-11=\",\"84-12=\",\"10+69=\",\"28-3=\",\"43+36=\")\n$newValues = @(\"1+82=\",\"9-4=\",\"20+63=\",\"7+8=\",\"18+33=\",\"90-68=\",\"21+46=\",\"31-2=\",\"74-20=\",\"60-13=\",\"70-22=\",\"71+22=\",\"61-40=\",\"7+44=\",\"87-62=\",\"23+7=\",\"74-7=\",\"2+78=\",\"49+33=\",\"61+35=\",\"84-80=\",\"28+6=\",\"44+26=\",\"45+20=\",\"98-19=\",\"62-60=\",\"13+63=\",\"17+70=\",\"46-28=\",\"55-50=\",\"41-31=\",\"75-60=\",\"82-59=\",\"29+59=\",\"1+84=\",\"38+17=\",\"61+7=\",\"99-12=\",\"39+53=\",\"35-20=\",\"53-6=\",\"41+6=\",\"94-10=\",\"94+3=\",\"20-19=\",\"20+52=\",\"37+31=\",\"78-24=\",\"99-2=\",\"58+0=\",\"43-20=\",\"83-3=\",\"66+8=\",\"72-36=\",\"6+67=\",\"28+25=\",\"58+30=\",\"59-20=\",\"26+15=\",\"5+75=\",\"94-82=\",\"57+42=\",\"41-21=\",\"86-70=\",\"71-63=\",\"69-52=\",\"96-34=\",\"93-30=\",\"17+60=\",\"37-13=\",\"30+52=\",\"51-1=\",\"74+21=\",\"24-2=\",\"51-40=\",\"85-27=\",\"66-20=\",\"23+56=\",\"7+84=\",\"40-0=\",\"95-83=\",\"48-8=\",\"28+59=\",\"38+6=\",\"54-3=\",\"79+1=\",\"79-77=\",\"30+14=\",\"97-65=\",\"66-54=\",\"54-10=\",\"41-7=\",\"85-57=\",\"38-31=\",\"72-69=\",\"59-5=\",\"51+0=\",\"7+47=\",\"68-7=\",\"36-23=\")\n\n# 1) Title paragraph\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -ne $titleOld) {\n    throw \"Title text mismatch: expected '$titleOld' got '$titleText'\"\n}\n$titlePara.Range.Text = $titleNew\n\n# 2) Table cells, row-major order\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\nif (($rows * $cols) -ne $oldValues.Count) {\n    throw \"Unexpected table size: $rows x $cols\"\n}\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $oldValues[$k]\n        $newVal = $newValues[$k]\n        if ($cellText -ne $expectedOld) {\n            throw \"Cell ($r,$c) text mismatch: expected '$expectedOld' got '$cellText'\"\n        }\n        if ($cellText -ne $newVal) {\n            $cell.Range.Text = $newVal\n        }\n        $k++\n    }\n}\n"}
